$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped crypto Price (column D) and Volume(1h) (column E)
# figures for this run. Both columns are stored as plain text in this
# sheet (coinranking scrape writes numbers/percentages as strings), so
# each new value is entered with a leading apostrophe to force Excel to
# keep it as text instead of reinterpreting it as a number/percentage.

$ws.Range("D2").Value = "'317.50"
$ws.Range("E2").Value = "'4.38%"
$ws.Range("E3").Value = "'10.95%"
$ws.Range("D4").Value = "'5.279"
$ws.Range("E4").Value = "'3.99%"
$ws.Range("E5").Value = "'3.29%"
$ws.Range("D6").Value = "'4.598"
$ws.Range("E6").Value = "'4.24%"
$ws.Range("D7").Value = "'1.339"
$ws.Range("E7").Value = "'33.38%"
$ws.Range("D8").Value = "'1.638"
$ws.Range("E8").Value = "'1.35%"
$ws.Range("D9").Value = "'0.1274"
$ws.Range("E9").Value = "'1.66%"
$ws.Range("D10").Value = "'0.1956"
$ws.Range("E10").Value = "'5.41%"
$ws.Range("D11").Value = "'0.09338"
$ws.Range("E11").Value = "'2.31%"
$ws.Range("D12").Value = "'0.04640"
$ws.Range("E12").Value = "'11.36%"
$ws.Range("D13").Value = "'0.1048"
$ws.Range("E13").Value = "'0.16%"
$ws.Range("D14").Value = "'0.001325"
$ws.Range("E14").Value = "'4.49%"
$ws.Range("D15").Value = "'0.04170"
$ws.Range("E15").Value = "'0.07%"
$ws.Range("D16").Value = "'0.005871"
$ws.Range("E16").Value = "'2.19%"
$ws.Range("D17").Value = "'3.338"
$ws.Range("E17").Value = "'0.24%"
$ws.Range("D18").Value = "'2.425"
$ws.Range("E18").Value = "'2.93%"
$ws.Range("E19").Value = "'4.55%"
$ws.Range("E20").Value = "'-4.45%"
$ws.Range("D21").Value = "'0.1385"
$ws.Range("E21").Value = "'-0.84%"
$ws.Range("D23").Value = "'0.001320"
$ws.Range("E23").Value = "'2.83%"
$ws.Range("D24").Value = "'0.004247"
$ws.Range("E24").Value = "'-5.52%"
$ws.Range("D25").Value = "'0.0001353"
$ws.Range("E25").Value = "'0.43%"
$ws.Range("D26").Value = "'0.0003547"
$ws.Range("E26").Value = "'-95.23%"
$ws.Range("D38").Value = "'0.02654"
$ws.Range("E38").Value = "'8.05%"
$ws.Range("D39").Value = "'0.05677"
$ws.Range("E39").Value = "'7.22%"
$ws.Range("D40").Value = "'0.01078"
$ws.Range("E40").Value = "'80.87%"
$ws.Range("D41").Value = "'0.008026"
$ws.Range("E41").Value = "'4.99%"
$ws.Range("D42").Value = "'0.1434"
$ws.Range("E42").Value = "'6.49%"
$ws.Range("D43").Value = "'0.007476"
$ws.Range("E43").Value = "'1.57%"
$ws.Range("D44").Value = "'0.008480"
$ws.Range("E44").Value = "'12.30%"
$ws.Range("D45").Value = "'0.3160"
$ws.Range("E45").Value = "'4.41%"
$ws.Range("D46").Value = "'0.00006645"
$ws.Range("E46").Value = "'-0.90%"
$ws.Range("E47").Value = "'0.40%"
$ws.Range("D48").Value = "'0.05489"
$ws.Range("E48").Value = "'31.53%"
$ws.Range("D49").Value = "'0.004008"
$ws.Range("E49").Value = "'-4.58%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.40%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.40%"
